$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7641556262969971
$ws.Range("B1").Value = 1.20589292049408
$ws.Range("C1").Value = 4.344531536102295
$ws.Range("D1").Value = 3.97257924079895
$ws.Range("E1").Value = 1.583240628242493
